$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.173.16"
$ws.Range("E2").Value = "  -4.11%  "
$ws.Range("D3").Value = "'1.655.29"
$ws.Range("E3").Value = "  -3.32%  "
$ws.Range("E4").Value = "  +0.31%  "
$ws.Range("D5").Value = "'216.37"
$ws.Range("E5").Value = "  -3.65%  "
$ws.Range("D6").Value = "'0.5130"
$ws.Range("E6").Value = "  -2.54%  "
$ws.Range("E7").Value = "  +0.24%  "
$ws.Range("D8").Value = "'0.2596"
$ws.Range("E8").Value = "  -1.68%  "
$ws.Range("D9").Value = "'0.06445"
$ws.Range("E9").Value = "  -3.20%  "
$ws.Range("D10").Value = "'19.75"
$ws.Range("E10").Value = "  -4.69%  "
$ws.Range("D11").Value = "'0.07815"
$ws.Range("E11").Value = "  +0.88%  "
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "'1.657.95"
$ws.Range("E12").Value = "  -3.26%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "'4.290"
$ws.Range("E13").Value = "  -3.76%  "
$ws.Range("D14").Value = "'1.886.09"
$ws.Range("E14").Value = "  -3.20%  "
$ws.Range("D15").Value = "'0.5513"
$ws.Range("E15").Value = "  -4.70%  "
$ws.Range("D16").Value = "0.0₅8003"
$ws.Range("D17").Value = "'64.08"
$ws.Range("E17").Value = "  -5.22%  "
$ws.Range("D18").Value = "'26.206.14"
$ws.Range("D19").Value = "'1.010"
$ws.Range("E19").Value = "  +0.34%  "
$ws.Range("D20").Value = "'208.59"
$ws.Range("E20").Value = "  -4.67%  "
$ws.Range("D21").Value = "'4.402"
$ws.Range("E21").Value = "  -5.06%  "
$ws.Range("D22").Value = "'10.08"
$ws.Range("E22").Value = "  -2.96%  "
$ws.Range("D23").Value = "'6.042"
$ws.Range("E23").Value = "  +0.41%  "
$ws.Range("D24").Value = "'1.011"
$ws.Range("E24").Value = "  +0.26%  "
$ws.Range("D25").Value = "'1.805"
$ws.Range("E25").Value = "  +5.80%  "
$ws.Range("D26").Value = "'144.49"
$ws.Range("E26").Value = "  -0.61%  "
$ws.Range("D27").Value = "'0.1173"
$ws.Range("E27").Value = "  -2.53%  "
$ws.Range("D28").Value = "'6.973"
$ws.Range("E28").Value = "  -3.39%  "
$ws.Range("D29").Value = "'15.81"
$ws.Range("E29").Value = "  -2.12%  "
$ws.Range("D30").Value = "'0.05068"
$ws.Range("E30").Value = "  -5.45%  "
$ws.Range("D31").Value = "'1.245"
$ws.Range("E31").Value = "  -3.82%  "
$ws.Range("D32").Value = "'3.352"
$ws.Range("E32").Value = "  -3.48%  "
$ws.Range("D33").Value = "'3.241"
$ws.Range("E33").Value = "  -3.97%  "
$ws.Range("D34").Value = "'1.552"
$ws.Range("E34").Value = "  -4.81%  "
$ws.Range("D35").Value = "'2.739"
$ws.Range("E35").Value = "  -3.82%  "
$ws.Range("D36").Value = "'2.358"
$ws.Range("E36").Value = "  -1.71%  "
$ws.Range("D37").Value = "'0.9206"
$ws.Range("E37").Value = "  -3.00%  "
$ws.Range("D38").Value = "'1.174.51"
$ws.Range("E38").Value = "  +1.70%  "
$ws.Range("D39").Value = "'0.5702"
$ws.Range("E39").Value = "  -2.81%  "
$ws.Range("D40").Value = "'0.01587"
$ws.Range("E40").Value = "  -3.71%  "
$ws.Range("D41").Value = "'2.569"
$ws.Range("E41").Value = "  -0.46%  "
$ws.Range("D42").Value = "'1.010"
$ws.Range("E42").Value = "  +0.25%  "
$ws.Range("D43").Value = "'5.662"
$ws.Range("E43").Value = "  -2.76%  "
$ws.Range("D44").Value = "'0.8259"
$ws.Range("E44").Value = "  -1.43%  "
$ws.Range("D45").Value = "'100.43"
$ws.Range("E45").Value = "  -0.50%  "
$ws.Range("D46").Value = "'1.799.37"
$ws.Range("E46").Value = "  -3.02%  "
$ws.Range("D47").Value = "0.0₈113"
$ws.Range("E47").Value = "  -3.52%  "
$ws.Range("D48").Value = "'0.4555"
$ws.Range("E48").Value = "  +0.04%  "
$ws.Range("D49").Value = "'55.39"
$ws.Range("E49").Value = "  -3.48%  "
$ws.Range("D50").Value = "'1.005"
$ws.Range("E50").Value = "  +0.07%  "
$ws.Range("D51").Value = "'7.863"
$ws.Range("E51").Value = "  -3.43%  "
